$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 361, shifting the existing row 361 (and all
# rows below it) down by one. This mirrors a new daily price record being
# added to the historical series.
$ws.Rows.Item(361).Insert()

# Populate the newly inserted row 361 with values. Most columns share the
# same constant values used throughout the sheet; only the date (D),
# volume (J), min/max/avg price (K/L/M) and price per kg (P) differ for
# this particular record.
$ws.Range("A361").Value = 4
$ws.Range("B361").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C361").Value = "Los Lagos"
$ws.Range("D361").Value = 45275
$ws.Range("E361").Value = 10
$ws.Range("F361").Value = 100112039
$ws.Range("G361").Value = "Ciboulette"
$ws.Range("H361").Value = "Sin especificar"
$ws.Range("I361").Value = "Primera"
$ws.Range("J361").Value = 240
$ws.Range("K361").Value = 2500
$ws.Range("L361").Value = 2500
$ws.Range("M361").Value = 2500
$ws.Range("N361").Value = "`$/docena de atados"
$ws.Range("O361").Value = "Región Metropolitana"
$ws.Range("P361").Value = 833
$ws.Range("Q361").Value = 3
$ws.Range("R361").Value = "Hortaliza"
